$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.934.54"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.107.12"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "576.79"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "178.04"
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.106.65"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "6.36"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").Value = "36.22"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "3.623.47"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "66.916.53"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "7.05"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "3.120.60"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "16.82"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "481.40"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("D22").Value = "7.81"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "0.691"
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("D24").Value = "83.73"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "12.59"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  -4.01%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "7.92"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").Value = "28.04"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("D34").Value = "0.0₃0941"
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "48.56"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("D37").Value = "5.61"
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("D38").Value = "0.942"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("D40").Value = "49.07"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").Value = "2.01"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "2.69"
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("D45").Value = "2.802.32"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "374.13"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0344"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "135.64"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D50").Value = "25.70"
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  +2.18%  "
